$d = $word.ActiveDocument

# The document currently starts with a single (empty) paragraph that
# carries the "_GoBack" bookmark. We need to insert four brand-new
# paragraphs *before* it, leaving that paragraph completely untouched.
#
# Using a zero-length (collapsed) Range at the very start of the story
# means InsertXML() will push the supplied paragraphs in front of the
# existing content without consuming/merging the original paragraph.

$insertionPoint = $d.Range(0, 0)

$bodyFragment = '<w:p>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Klark</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> gua </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>underground</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> viva el </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>atrtr</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>' +
    '<w:p>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Kjfshvlfh</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>' +
    '<w:p>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Klark</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>' +
    '<w:p>' +
    '<w:r><w:t xml:space="preserve">Juan </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>carlos</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> cabrera </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>roman</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $bodyFragment + '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData>' +
    '</pkg:part>' +
    '</pkg:package>'

$insertionPoint.InsertXML($xml)
